$d = $word.ActiveDocument

# 1. Add w:bCs to the bold run in the "Please note..." paragraph.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Please note that the steps show rounded numbers*") {
        foreach ($r in $p.Range.Words) {
        }
        $p.Range.Font.BoldBi = 1
    }
}
